# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) for the second data row
# (the 44e77bcf... file) on both the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-25 09:21:14"
$wsZhCn.Range("H3").Value = "2016-03-25 09:22:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-25 09:21:19"
$wsDeDe.Range("H3").Value = "2016-03-25 09:22:22"
